$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to stay text while we overwrite their values,
# so Excel does not auto-convert numeric-looking strings (e.g. "1.008") to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.371.45'
$ws.Range("E2").Value = '  -3.67%  '
$ws.Range("D3").Value = '1.665.00'
$ws.Range("E3").Value = '  -2.57%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = '218.61'
$ws.Range("E5").Value = '  -2.33%  '
$ws.Range("D6").Value = '0.5151'
$ws.Range("E6").Value = '  -3.10%  '
$ws.Range("D7").Value = '1.007'
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("D8").Value = '0.06439'
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("D9").Value = '0.2570'
$ws.Range("E9").Value = '  -3.24%  '
$ws.Range("D10").Value = '19.92'
$ws.Range("E10").Value = '  -4.29%  '
$ws.Range("D11").Value = '0.07668'
$ws.Range("E11").Value = '  +0.32%  '
$ws.Range("D12").Value = '4.335'
$ws.Range("E12").Value = '  -5.17%  '
$ws.Range("D13").Value = '1.664.14'
$ws.Range("E13").Value = '  -2.81%  '
$ws.Range("D14").Value = '1.894.43'
$ws.Range("E14").Value = '  -2.57%  '
$ws.Range("D15").Value = '0.5529'
$ws.Range("E15").Value = '  -3.35%  '
$ws.Range("D16").Value = '0.0₅8024'
$ws.Range("E16").Value = '  -1.69%  '
$ws.Range("D17").Value = '64.59'
$ws.Range("E17").Value = '  -4.82%  '
$ws.Range("D18").Value = '26.401.15'
$ws.Range("E18").Value = '  -3.49%  '
$ws.Range("D19").Value = '1.007'
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("D20").Value = '209.52'
$ws.Range("E20").Value = '  -3.16%  '
$ws.Range("D21").Value = '4.411'
$ws.Range("E21").Value = '  -5.51%  '
$ws.Range("D22").Value = '10.08'
$ws.Range("E22").Value = '  -3.25%  '
$ws.Range("D23").Value = '5.869'
$ws.Range("E23").Value = '  -1.69%  '
$ws.Range("D25").Value = '145.41'
$ws.Range("E25").Value = '  +2.55%  '
$ws.Range("D26").Value = '1.734'
$ws.Range("E26").Value = '  -2.00%  '
$ws.Range("D28").Value = '6.988'
$ws.Range("E28").Value = '  -3.88%  '
$ws.Range("D29").Value = '15.76'
$ws.Range("E29").Value = '  -3.33%  '
$ws.Range("D30").Value = '0.05224'
$ws.Range("E30").Value = '  -3.46%  '
$ws.Range("D31").Value = '1.261'
$ws.Range("E31").Value = '  -2.35%  '
$ws.Range("D32").Value = '3.365'
$ws.Range("E32").Value = '  -3.99%  '
$ws.Range("D33").Value = '3.212'
$ws.Range("E33").Value = '  -6.26%  '
$ws.Range("D34").Value = '1.570'
$ws.Range("E34").Value = '  -4.45%  '
$ws.Range("D35").Value = '2.753'
$ws.Range("E35").Value = '  -4.29%  '
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("D37").Value = '0.9249'
$ws.Range("E37").Value = '  -2.58%  '
$ws.Range("D38").Value = '0.5707'
$ws.Range("E38").Value = '  -2.61%  '
$ws.Range("D39").Value = '1.151.40'
$ws.Range("E39").Value = '  +10.11%  '
$ws.Range("D40").Value = '0.01595'
$ws.Range("E40").Value = '  -2.07%  '
$ws.Range("D41").Value = '1.007'
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("D42").Value = '0.8406'
$ws.Range("E42").Value = '  -0.43%  '
$ws.Range("D43").Value = '5.635'
$ws.Range("E43").Value = '  -3.89%  '
$ws.Range("D44").Value = '100.16'
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("D45").Value = '1.804.94'
$ws.Range("E45").Value = '  -2.56%  '
$ws.Range("D46").Value = '0.0₈111'
$ws.Range("E46").Value = '  -6.30%  '
$ws.Range("D47").Value = '0.4494'
$ws.Range("E47").Value = '  -0.33%  '
$ws.Range("D48").Value = '55.83'
$ws.Range("E48").Value = '  -3.68%  '
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").Value = '7.935'
$ws.Range("E50").Value = '  -1.79%  '
$ws.Range("D51").Value = '0.05112'
$ws.Range("E51").Value = '  -2.55%  '

# Restore the original (default/general) cell style now that the text values are locked in.
$ws.Range("D2:E51").Style = "Normal"
